# Add two new rows (16 and 17) of support-ticket data to Sheet1, just
# below the existing data (which ends at row 15), extending the used
# range from A1:C15 to A1:C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 ---------------------------------------------------------
# Column A: alphanumeric "phone-like" string -> plain text already.
$ws.Range("A16").Value = "0777553527236346346523724y635u7"

# Column B: purely numeric-looking string with a leading zero, so it
# must be forced to text (else Excel would coerce it to a Number and
# drop the leading zero). Apply a text format before assigning, then
# restore the cell style to Normal so no stray formatting is left
# behind on the cell.
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "0743555526"
$ws.Range("B16").Style = "Normal"

# Column C: empty issue field for this row.
$ws.Range("C16").Value = ""

# --- Row 17 ---------------------------------------------------------
$ws.Range("A17").Value = "supun dissanayaka"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "0382250162"
$ws.Range("B17").Style = "Normal"

$ws.Range("C17").Value = "my router is not working properly there is red light blinking on los bulb"
